# Bento object repository revisited
#
# The "FilesTab" Cypher query (cell B4 on the "startup" sheet) was
# simplified: the `File Type` and `Breed` projections were dropped from
# the RETURN clause. Re-write the cell with the updated query text and
# move the selection onto it, matching the authored edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

$query = @"
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
 MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
 MATCH (samp:sample)-->(c) 
 WHERE samp.specific_sample_pathology IN ["T Cell Lymphoma"]  
WITH DISTINCT f, parent, c, demo, diag, s
RETURN coalesce(f.file_name, '') AS ``File Name``, 
        coalesce(labels(parent)[0], '') AS ``Association``,
        coalesce(f.file_description, '') AS ``Description``,
        coalesce(f.file_format, '') AS ``Format``,
        coalesce(f.file_size, '') AS ``Size``,
        coalesce(c.case_id, '') AS ``Case ID``, 
        coalesce(diag.disease_term,'') AS Diagnosis , 
        coalesce(s.clinical_study_designation,'') AS ``Study Code``
"@

$cell = $ws.Range("B4")
$cell.Value = $query

# Reflect the author's updated cursor/selection position on the sheet.
$cell.Select()
